$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (row 2) is being re-shuffled/extended to match the regenerated
# staging template schema: a new "BusinessKey" column becomes the first
# column, the old headers shift right, and a new 4th column is appended
# with what used to be the 3rd column's heading.
$ws.Range("A2").Value = "BusinessKey"
$ws.Range("B2").Value = "Framework_Indicator_ID"
$ws.Range("C2").Value = "FrameworkBusinessKey"

# New column D2, matching the bold/underline header style used by A2:C2
$ws.Range("D2").Value = "IndicatorBusinessKey"
$ws.Range("D2").Font.Bold = $true
$ws.Range("D2").Font.Underline = $true
